$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "As a technical writer for a software company, you are creating a Word document detailing its features. You have listed all the features by pressing the Enter key after each item. However, you feel this content would be easier to read as a bulleted list.How should you achieve this most efficiently?",
        "ques_type": 2,
        "options": [
            "Click the bulleted list icon in the Home toolbar &gt Move the items into the new list using the Backspace key.",
            "Highlight all the items in the list &gt Click the bulleted list icon in the Home toolbar.",
            "Highlight the first item in the list &gt Click the bulleted list icon in the Home toolbar.",
            "Click the bulleted list icon in the Home toolbar &gt Highlight all the items in the list."
        ],
        "score": "Highlight all the items in the list &gt Click the bulleted list icon in the Home toolbar."
    },
    {
        "title": "As a data analyst in an educational institution, you've been given a large Excel spreadsheet of student names along with their respective GPA scores. The Dean has asked you to identify the student with the highest GPA score.Which tool should you use to achieve this most efficiently?",
        "ques_type": 2,
        "options": [
            "Filter",
            "Sort",
            "Find",
            "Search"
        ],
        "score": "Sort"
    },
    {
        "title": "You are going to present PowerPoint slides during an upcoming sales meeting. You plan to connect your laptop to the projector and describe the product to the audience. To not miss anything, you want to have talking points near your eyes when you give your presentation.Which actions should you take to present most professionally?",
        "ques_type": 15,
        "options": [
            "Write your talking points in a separate Word document. ",
            "Open the Review tab &gt add notes to the slides as comments.",
            "Open the View tab &gt enable the Notes Pane &gt add notes to the slides.",
            "Print the PowerPoint file using the Notes Pages layout.",
            "Refer to a separate Word document during the presentation.",
            "Use Presenter View during the presentation."
        ],
        "score": [
            "Open the View tab &gt enable the Notes Pane &gt add notes to the slides.",
            "Use Presenter View during the presentation."
        ]
    },
    {
        "title": "You are an administrative assistant who uses Outlook for email management. Your manager frequently sends you important updates and tasks via email, which often get lost among other messages. Your goal is to make these messages more accessible.What action should you take to achieve it most efficiently?",
        "ques_type": 2,
        "options": [
            "Create a rule to mark all manager\u2019s emails as important.",
            "Create a rule to categorize manager\u2019s emails automatically.",
            "Set up a customized notification for the manager's emails.",
            "Create a separate folder and set up a rule to move manager\u2019s emails to this folder."
        ],
        "score": "Create a separate folder and set up a rule to move manager\u2019s emails to this folder."
    }
]
'@

$ws.Rows(1).Delete()
$ws.Range("A1").Value = $newText
$ws.Rows(1).AutoFit()
